$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.890.35'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.373.98'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').Value = '2.370.84'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '2.794.32'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('D17').Value = '59.772.28'
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('D18').Value = '2.366.09'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.93%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '63.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '556.33'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.24%  '
$ws.Range('D30').Value = '0.0₃0924'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('E32').Value = '  -3.22%  '
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('E36').Value = '  +3.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '153.43'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.13'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '139.73'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0500'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('E51').Value = '  -1.43%  '
